# Append one new data-log row to the bottom of each of the four
# worksheets (ROW50-FE-LIFTER, ROW50-MID-LIFTER, ROW11-FE-LIFTER,
# ROW11-MID-LIFTER), matching the layout/format of the existing rows.

$wb = $excel.ActiveWorkbook

# Large ID_DEC value shared by most of the new rows (same constant that
# already appears throughout the sheets). Built via a string->double cast
# because the interpreter's numeric literal parser does not understand
# scientific notation (e.g. "5.68631262647114e+23") directly.
$bigId = [double]"5.68631262647114e+23"

# ---------------------------------------------------------------------
# Sheet 1: ROW50-FE-LIFTER -> append row 78
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$r = 78
$p = $r - 1
$ws1.Range("A$r").NumberFormat = $ws1.Range("A$p").NumberFormat
$ws1.Range("A$r").Value = 45763.26375638889
$ws1.Range("B$r").Value = "0x01,0x90"
$ws1.Range("C$r").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D$r").Value = "0x01,0x4a"
$ws1.Range("E$r").Value = "0xe"
$ws1.Range("F$r").Value = 400
$ws1.Range("G$r").Value = $bigId
$ws1.Range("H$r").Value = 330
$ws1.Range("I$r").Value = 14

# ---------------------------------------------------------------------
# Sheet 2: ROW50-MID-LIFTER -> append row 80
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$r = 80
$p = $r - 1
$ws2.Range("A$r").NumberFormat = $ws2.Range("A$p").NumberFormat
$ws2.Range("A$r").Value = 45763.22726851852
$ws2.Range("B$r").Value = "0x01,0x90 "
$ws2.Range("C$r").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D$r").Value = "0x01,0x4e"
$ws2.Range("E$r").Value = "0x19"
$ws2.Range("F$r").Value = 400
# This sheet stores the ID_DEC column as plain text (matches every other
# row in this particular sheet), so force a text value rather than
# letting it coerce into a numeric cell. Reset the style back to Normal
# afterwards so only the cell's data type (not its number format) is
# affected, matching the rest of the column.
$ws2.Range("G$r").NumberFormat = "@"
$ws2.Range("G$r").Value = "568631262647113771663628"
$ws2.Range("G$r").Style = "Normal"
$ws2.Range("H$r").Value = 334
$ws2.Range("I$r").Value = 25

# ---------------------------------------------------------------------
# Sheet 3: ROW11-FE-LIFTER -> append row 78
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$r = 78
$p = $r - 1
$ws3.Range("A$r").NumberFormat = $ws3.Range("A$p").NumberFormat
$ws3.Range("A$r").Value = 45763.297658125
$ws3.Range("B$r").Value = "0x01,0x90"
$ws3.Range("C$r").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D$r").Value = "0x01,0x4a"
$ws3.Range("E$r").Value = "0x14"
$ws3.Range("F$r").Value = 400
$ws3.Range("G$r").Value = $bigId
$ws3.Range("H$r").Value = 330
$ws3.Range("I$r").Value = 20

# ---------------------------------------------------------------------
# Sheet 4: ROW11-MID-LIFTER -> append row 78
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$r = 78
$p = $r - 1
$ws4.Range("A$r").NumberFormat = $ws4.Range("A$p").NumberFormat
$ws4.Range("A$r").Value = 45763.42206883102
$ws4.Range("B$r").Value = "0x01,0x90"
$ws4.Range("C$r").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D$r").Value = "0x01,0x52"
$ws4.Range("E$r").Value = "0x19"
$ws4.Range("F$r").Value = 400
$ws4.Range("G$r").Value = $bigId
$ws4.Range("H$r").Value = 338
$ws4.Range("I$r").Value = 25
